$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (was 강원춘천시막장FC / A=0) -> FC리틀슛 / A=2
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "FC리틀슛"
$ws.Range("C4").Value = "창단일자`n20070301`n주소`n부산 북구 화명대로94번길 45 (화명동, 화명그린힐아파트) 103동105호`n연락처`n010-4701-8300`n감독`n박성용`n팀 관리에서 팀 SNS를 등록하세요"

# Row 5 (was CRASSODONFC / A=0) -> FC목포 / A=3
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "FC목포"
$ws.Range("C5").Value = "창단일자`n20091224`n주소`n전남 목포시 내화마을길 89 (대양동, 목포국제축구센터)`n연락처`n061-274-0171`n감독`n조덕제`n팀 관리에서 팀 SNS를 등록하세요"

# Row 6 (was 강원춘천시막장FC / A=0) -> FC서울 / A=4
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "FC서울"
$ws.Range("C6").Value = "창단일자`n19831222`n주소`n서울특별시 마포구 월드컵로 240 (성산동, 월드컵주경기장) 서울월드컵경기장 서측 4층 FC서울 사무실`n연락처`n02-376-3946`n감독`n안익수`n팀 관리에서 팀 SNS를 등록하세요"

# Row 7 (was CRASSODONFC / A=0) -> FC서울U18서울오산고 / A=5
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "FC서울U18서울오산고"
$ws.Range("C7").Value = "창단일자`n20121227`n주소`n서울특별시 마포구 월드컵로 240 (성산동, FC서울) FC서울 운영육성팀`n연락처`n02-376-3044`n감독`n윤현필`n팀 관리에서 팀 SNS를 등록하세요"

# Row 8 (new) -> FC안양 / A=6
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "FC안양"
$ws.Range("C8").Value = "창단일자`n20130202`n주소`n경기도 안양시 동안구 평촌대로 389 (비산동, 안양종합운동장) 안양종합운동장 내 FC안양 사무국`n연락처`n031-476-3500`n감독`n이우형`n팀 관리에서 팀 SNS를 등록하세요"

# Row 9 (new) -> FC증산 / A=7
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "FC증산"
$ws.Range("C9").Value = "창단일자`n19970917`n주소`n경남 양산시 물금읍 물금로 9 (물금리) 더 스퀘어 201호`n연락처`n010-4635-5798`n감독`n김효준1`n팀 관리에서 팀 SNS를 등록하세요"

# Row 10 (new) -> FC안양 / A=8
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "FC안양"
$ws.Range("C10").Value = "창단일자`n20130202`n주소`n경기도 안양시 동안구 평촌대로 389 (비산동, 안양종합운동장) 안양종합운동장 내 FC안양 사무국`n연락처`n031-476-3500`n감독`n이우형`n팀 관리에서 팀 SNS를 등록하세요"

# Row 11 (new) -> LS축구아카데미 / A=9
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "LS축구아카데미"
$ws.Range("C11").Value = "창단일자`n20200805`n주소`n대구 달서구 상화북로 189 (상인동) 신동빌딩4층`n연락처`n010-8608-3439`n감독`n이진표`n팀 관리에서 팀 SNS를 등록하세요"

# Row 12 (new) -> 강원고성군간성클럽 / A=0
$ws.Range("A12").Value = 0
$ws.Range("B12").Value = "강원고성군간성클럽"
$ws.Range("C12").Value = "창단일자`n20140302`n주소`n강원 고성군 간성읍 수성로 111 (상리, 종합운동장)`n연락처`n010-2171-6033`n감독`n박득쇠`n팀 관리에서 팀 SNS를 등록하세요"

# Row 13 (new) -> 강원고성군거진FC / A=0
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = "강원고성군거진FC"
$ws.Range("C13").Value = "창단일자`n20210101`n주소`n강원도 고성군 거진읍 거진운동장길 30 (거진리) 1`n연락처`n010-4677-7077`n감독`n팀 관리에서 팀 SNS를 등록하세요"

# Row 14 (new) -> 강원고성군고성사랑축구회 / A=0
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = "강원고성군고성사랑축구회"
$ws.Range("C14").Value = "창단일자`n20180102`n주소`n강원 고성군 토성면 장새미1길 11 (용암리) 용암리`n연락처`n010-2828-4338`n감독`n팀 관리에서 팀 SNS를 등록하세요"

# Copy the header-row cell formatting (bold font, thin border, center/top alignment)
# from the existing A4 (style index 1) into the newly added A8:A14 cells so they
# match the style used by the other column-A cells, without introducing new styles.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A8:A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Writing multi-line strings into column C nudges the engine into giving those
# rows an explicit custom row height; put the affected rows back to the
# sheet's normal auto-sized height so no stray row-height formatting is left
# behind beyond what the source data describes.
$ws.Range("C4:C14").EntireRow.AutoFit()
